# Auto-generated edit script applying scheduled-runner price refresh to Sheets
# Each sheet corresponds to a crafting job (Table_<JOB>); columns H-N hold
# market-price-derived figures that get refreshed by the scheduled runner.
$wb = $excel.ActiveWorkbook

# ----- Sheet: ALC -----
$ws = $wb.Worksheets.Item("ALC")
# Row 28
$ws.Cells.Item(28, 8).Value = 175.72728
$ws.Cells.Item(28, 9).Value = 173.3
$ws.Cells.Item(28, 11).Value = 173.3
$ws.Cells.Item(28, 13).Value = 311.7
# Row 70
$ws.Cells.Item(70, 8).Value = 3031.8845
$ws.Cells.Item(70, 9).Value = 1158.2
$ws.Cells.Item(70, 10).Value = 3478
$ws.Cells.Item(70, 11).Value = 3474.6
$ws.Cells.Item(70, 12).Value = 10434
$ws.Cells.Item(70, 13).Value = -3204.6
$ws.Cells.Item(70, 14).Value = -10974
# Row 73
$ws.Cells.Item(73, 8).Value = 3031.8845
$ws.Cells.Item(73, 9).Value = 1158.2
$ws.Cells.Item(73, 10).Value = 3478
$ws.Cells.Item(73, 11).Value = 3474.6
$ws.Cells.Item(73, 12).Value = 10434
$ws.Cells.Item(73, 13).Value = -2538.6
$ws.Cells.Item(73, 14).Value = -12306

# ----- Sheet: ARM -----
$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Cells.Item(2, 8).Value = 5954620
$ws.Cells.Item(2, 9).Value = 14707228
$ws.Cells.Item(2, 10).Value = 2847.12
$ws.Cells.Item(2, 11).Value = 14707228
$ws.Cells.Item(2, 12).Value = 2847.12
$ws.Cells.Item(2, 13).Value = -14707115
$ws.Cells.Item(2, 14).Value = -3073.12
# Row 23
$ws.Cells.Item(23, 8).Value = 24273.182
$ws.Cells.Item(23, 9).Value = 40002.5
$ws.Cells.Item(23, 10).Value = 20777.777
$ws.Cells.Item(23, 11).Value = 40002.5
$ws.Cells.Item(23, 12).Value = 20777.777
$ws.Cells.Item(23, 13).Value = -39743.5
$ws.Cells.Item(23, 14).Value = -21295.777
# Row 104
$ws.Cells.Item(104, 8).Value = 250000
$ws.Cells.Item(104, 10).Value = 250000
$ws.Cells.Item(104, 12).Value = 250000
$ws.Cells.Item(104, 14).Value = -256988
# Row 116
$ws.Cells.Item(116, 8).Value = 5954620
$ws.Cells.Item(116, 9).Value = 14707228
$ws.Cells.Item(116, 10).Value = 2847.12
$ws.Cells.Item(116, 11).Value = 14707228
$ws.Cells.Item(116, 12).Value = 2847.12
$ws.Cells.Item(116, 13).Value = -14704934
$ws.Cells.Item(116, 14).Value = -7435.12

# ----- Sheet: BSM -----
$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Cells.Item(3, 8).Value = 5954620
$ws.Cells.Item(3, 9).Value = 14707228
$ws.Cells.Item(3, 10).Value = 2847.12
$ws.Cells.Item(3, 11).Value = 14707228
$ws.Cells.Item(3, 12).Value = 2847.12
$ws.Cells.Item(3, 13).Value = -14707114
$ws.Cells.Item(3, 14).Value = -3075.12
# Row 20
$ws.Cells.Item(20, 8).Value = 1341.2258
$ws.Cells.Item(20, 9).Value = 821.8889
$ws.Cells.Item(20, 10).Value = 2060.3076
$ws.Cells.Item(20, 11).Value = 821.8889
$ws.Cells.Item(20, 12).Value = 2060.3076
$ws.Cells.Item(20, 13).Value = -574.8889
$ws.Cells.Item(20, 14).Value = -2554.3076
# Row 94
$ws.Cells.Item(94, 8).Value = 571.26666
$ws.Cells.Item(94, 9).Value = 552.63635
$ws.Cells.Item(94, 10).Value = 622.5
$ws.Cells.Item(94, 11).Value = 552.63635
$ws.Cells.Item(94, 12).Value = 622.5
$ws.Cells.Item(94, 13).Value = -101.63635
$ws.Cells.Item(94, 14).Value = -1524.5
# Row 134
$ws.Cells.Item(134, 8).Value = 2378.8667
$ws.Cells.Item(134, 9).Value = 1680.9546
$ws.Cells.Item(134, 10).Value = 4298.125
$ws.Cells.Item(134, 11).Value = 5042.8638
$ws.Cells.Item(134, 12).Value = 12894.375
$ws.Cells.Item(134, 13).Value = -2507.8638
$ws.Cells.Item(134, 14).Value = -17964.375

# ----- Sheet: CRP -----
$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Cells.Item(16, 8).Value = 2611.182
$ws.Cells.Item(16, 9).Value = 1422.2
$ws.Cells.Item(16, 10).Value = 3602
$ws.Cells.Item(16, 11).Value = 1422.2
$ws.Cells.Item(16, 12).Value = 3602
$ws.Cells.Item(16, 13).Value = -1135.2
$ws.Cells.Item(16, 14).Value = -4176
# Row 31
$ws.Cells.Item(31, 8).Value = 2782198.8
$ws.Cells.Item(31, 9).Value = 4547195.5
$ws.Cells.Item(31, 10).Value = 8632.071
$ws.Cells.Item(31, 11).Value = 4547195.5
$ws.Cells.Item(31, 12).Value = 8632.071
$ws.Cells.Item(31, 13).Value = -4546900.5
$ws.Cells.Item(31, 14).Value = -9222.071
# Row 34
$ws.Cells.Item(34, 8).Value = 2782198.8
$ws.Cells.Item(34, 9).Value = 4547195.5
$ws.Cells.Item(34, 10).Value = 8632.071
$ws.Cells.Item(34, 11).Value = 4547195.5
$ws.Cells.Item(34, 12).Value = 8632.071
$ws.Cells.Item(34, 13).Value = -4546993.5
$ws.Cells.Item(34, 14).Value = -9036.071
# Row 58
$ws.Cells.Item(58, 8).Value = 31254044
$ws.Cells.Item(58, 9).Value = 3119.8
$ws.Cells.Item(58, 10).Value = 45459010
$ws.Cells.Item(58, 11).Value = 3119.8
$ws.Cells.Item(58, 12).Value = 45459010
$ws.Cells.Item(58, 13).Value = -2916.8
$ws.Cells.Item(58, 14).Value = -45459416
# Row 92
$ws.Cells.Item(92, 8).Value = 0
$ws.Cells.Item(92, 10).Value = 0
$ws.Cells.Item(92, 14).ClearContents()
# Row 105
$ws.Cells.Item(105, 8).Value = 1986.0667
$ws.Cells.Item(105, 9).Value = 1343.6364
$ws.Cells.Item(105, 10).Value = 3752.75
$ws.Cells.Item(105, 11).Value = 1343.6364
$ws.Cells.Item(105, 12).Value = 3752.75
$ws.Cells.Item(105, 13).Value = 403.3635999999999
$ws.Cells.Item(105, 14).Value = -7246.75
# Row 113
$ws.Cells.Item(113, 8).Value = 2611.182
$ws.Cells.Item(113, 9).Value = 1422.2
$ws.Cells.Item(113, 10).Value = 3602
$ws.Cells.Item(113, 11).Value = 1422.2
$ws.Cells.Item(113, 12).Value = 3602
$ws.Cells.Item(113, 13).Value = 747.8
$ws.Cells.Item(113, 14).Value = -7942
# Row 136
$ws.Cells.Item(136, 8).Value = 31254044
$ws.Cells.Item(136, 9).Value = 3119.8
$ws.Cells.Item(136, 10).Value = 45459010
$ws.Cells.Item(136, 11).Value = 9359.400000000001
$ws.Cells.Item(136, 12).Value = 136377030
$ws.Cells.Item(136, 13).Value = -6809.400000000001
$ws.Cells.Item(136, 14).Value = -136382130

# ----- Sheet: CUL -----
$ws = $wb.Worksheets.Item("CUL")
# Row 63
$ws.Cells.Item(63, 8).Value = 3253
$ws.Cells.Item(63, 9).Value = 1012
$ws.Cells.Item(63, 10).Value = 4000
$ws.Cells.Item(63, 11).Value = 3036
$ws.Cells.Item(63, 12).Value = 12000
$ws.Cells.Item(63, 13).Value = -2287
$ws.Cells.Item(63, 14).Value = -13498
# Row 66
$ws.Cells.Item(66, 8).Value = 3253
$ws.Cells.Item(66, 9).Value = 1012
$ws.Cells.Item(66, 10).Value = 4000
$ws.Cells.Item(66, 11).Value = 9108
$ws.Cells.Item(66, 12).Value = 36000
$ws.Cells.Item(66, 13).Value = -5364
$ws.Cells.Item(66, 14).Value = -43488
# Row 69
$ws.Cells.Item(69, 8).Value = 93976.164
$ws.Cells.Item(69, 9).Value = 1450
$ws.Cells.Item(69, 11).Value = 4350
$ws.Cells.Item(69, 13).Value = -3539
# Row 72
$ws.Cells.Item(72, 8).Value = 93976.164
$ws.Cells.Item(72, 9).Value = 1450
$ws.Cells.Item(72, 11).Value = 13050
$ws.Cells.Item(72, 13).Value = -8994
# Row 92
$ws.Cells.Item(92, 8).Value = 1169.6818
$ws.Cells.Item(92, 9).Value = 790
$ws.Cells.Item(92, 10).Value = 1229.6316
$ws.Cells.Item(92, 11).Value = 2370
$ws.Cells.Item(92, 12).Value = 3688.8948
$ws.Cells.Item(92, 13).Value = -1122
$ws.Cells.Item(92, 14).Value = -6184.8948
# Row 113
$ws.Cells.Item(113, 8).Value = 729.1539
$ws.Cells.Item(113, 9).Value = 585.75
$ws.Cells.Item(113, 10).Value = 792.8889
$ws.Cells.Item(113, 11).Value = 1757.25
$ws.Cells.Item(113, 12).Value = 2378.6667
$ws.Cells.Item(113, 13).Value = 412.75
$ws.Cells.Item(113, 14).Value = -6718.6667

# ----- Sheet: GSM -----
$ws = $wb.Worksheets.Item("GSM")
# Row 132
$ws.Cells.Item(132, 8).Value = 3293.5
$ws.Cells.Item(132, 9).Value = 1966.4445
$ws.Cells.Item(132, 10).Value = 4999.7144
$ws.Cells.Item(132, 11).Value = 5899.333500000001
$ws.Cells.Item(132, 12).Value = 14999.1432
$ws.Cells.Item(132, 13).Value = -3369.333500000001
$ws.Cells.Item(132, 14).Value = -20059.1432

# ----- Sheet: LTW -----
$ws = $wb.Worksheets.Item("LTW")
# Row 16
$ws.Cells.Item(16, 8).Value = 9380
$ws.Cells.Item(16, 10).Value = 11000
$ws.Cells.Item(16, 12).Value = 11000
$ws.Cells.Item(16, 14).Value = -11340
# Row 46
$ws.Cells.Item(46, 8).Value = 2519.125
$ws.Cells.Item(46, 9).Value = 500
$ws.Cells.Item(46, 10).Value = 3192.1667
$ws.Cells.Item(46, 11).Value = 500
$ws.Cells.Item(46, 12).Value = 3192.1667
$ws.Cells.Item(46, 13).Value = -312
$ws.Cells.Item(46, 14).Value = -3568.1667
# Row 136
$ws.Cells.Item(136, 8).Value = 4548577
$ws.Cells.Item(136, 9).Value = 8335908
$ws.Cells.Item(136, 10).Value = 3780
$ws.Cells.Item(136, 11).Value = 25007724
$ws.Cells.Item(136, 12).Value = 11340
$ws.Cells.Item(136, 13).Value = -25005174
$ws.Cells.Item(136, 14).Value = -16440

# ----- Sheet: WVR -----
$ws = $wb.Worksheets.Item("WVR")
# Row 2
$ws.Cells.Item(2, 8).Value = 2598000.5
$ws.Cells.Item(2, 9).Value = 8999
$ws.Cells.Item(2, 10).Value = 3115800.8
$ws.Cells.Item(2, 11).Value = 8999
$ws.Cells.Item(2, 12).Value = 3115800.8
$ws.Cells.Item(2, 13).Value = -8887
$ws.Cells.Item(2, 14).Value = -3116024.8
# Row 81
$ws.Cells.Item(81, 8).Value = 749.8889
$ws.Cells.Item(81, 9).Value = 570
$ws.Cells.Item(81, 10).Value = 974.75
$ws.Cells.Item(81, 11).Value = 1140
$ws.Cells.Item(81, 12).Value = 1949.5
$ws.Cells.Item(81, 13).Value = -79
$ws.Cells.Item(81, 14).Value = -4071.5
# Row 84
$ws.Cells.Item(84, 8).Value = 749.8889
$ws.Cells.Item(84, 9).Value = 570
$ws.Cells.Item(84, 10).Value = 974.75
$ws.Cells.Item(84, 11).Value = 5700
$ws.Cells.Item(84, 12).Value = 9747.5
$ws.Cells.Item(84, 13).Value = -396
$ws.Cells.Item(84, 14).Value = -20355.5
# Row 132
$ws.Cells.Item(132, 8).Value = 275029.75
$ws.Cells.Item(132, 9).Value = 436144.9
$ws.Cells.Item(132, 10).Value = 10340.571
$ws.Cells.Item(132, 11).Value = 1308434.7
$ws.Cells.Item(132, 12).Value = 31021.713
$ws.Cells.Item(132, 13).Value = -1305904.7
$ws.Cells.Item(132, 14).Value = -36081.713

Write-Host "Applied 211 cell updates and 1 cell clear(s)."
